$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "excel (3)" to "excel (4)"
$ws.Name = "excel (4)"

# Rename the matching defined name from excel__3 to excel__4
# (RefersTo auto-follows the sheet rename)
$names = $wb.Names
foreach ($n in $names) {
    if ($n.Name -like "*excel__3*") {
        $n.Name = "excel__4"
    }
}

# Update the three count cells
$ws.Range("K17").Value2 = 2
$ws.Range("K19").Value2 = 6
$ws.Range("K21").Value2 = 9

# Update the hyperlink display text (moreBaseParts/moreNew/moreReused 6 -> 3)
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$13') {
        $h.TextToDisplay = "http://localhost:2468/reports/sizeest.class?moreBaseParts=3"
    }
    if ($addr -eq '$A$23') {
        $h.TextToDisplay = "http://localhost:2468/reports/sizeest.class?moreNew=3"
    }
    if ($addr -eq '$A$27') {
        $h.TextToDisplay = "http://localhost:2468/reports/sizeest.class?moreReused=3"
    }
}
